$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.239943265914917
$ws.Range("B1").Value = 2.359054327011108
$ws.Range("C1").Value = 3.775464534759521
$ws.Range("D1").Value = 3.280536651611328
$ws.Range("E1").Value = 1.258055090904236
